$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (y_0_forecast) and Column E (y_1_forecast) values per row
$values = @{
    4  = @{ C = 0.26989196962941;    E = 0.5668332592311476 }
    5  = @{ C = 3.872616460773104;   E = 2.302293339151018 }
    6  = @{ C = 1.385102141502959;   E = 1.418650881502459 }
    7  = @{ C = -0.09990319152052285; E = 0.9503446311848185 }
    8  = @{ C = 1.985365612881851;   E = 1.683270396159919 }
    9  = @{ C = 1.326941625882871;   E = 1.183026708734536 }
    10 = @{ C = 1.606236217798274;   E = 1.5164178751083 }
    11 = @{ C = 1.730343268967593;   E = 1.500360889159746 }
    12 = @{ C = 2.175463816693268;   E = 1.290682042302871 }
    13 = @{ C = 0.7385331577992593;  E = 1.193318741914795 }
    14 = @{ C = -1.986210268830169; E = -1.030455917249229 }
    15 = @{ C = 0.4846423081591666;  E = 1.007198277338284 }
    16 = @{ C = 1.990013243928312;   E = 1.178997008351645 }
    17 = @{ C = -0.07210020592836042; E = 0.7849798646630823 }
    18 = @{ C = -0.03189435474734159; E = 0.8280060478212947 }
    19 = @{ C = 0.1720146172997206;  E = 0.822899916144304 }
}

foreach ($row in $values.Keys) {
    $ws.Range("C$row").Value = $values[$row].C
    $ws.Range("E$row").Value = $values[$row].E
}
